# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stock) worksheet (5th sheet) gets three new trailing
# columns: date / legislator_name / legislator_id, populated for every
# existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

$legDate = "2012-04-18"
$legName = "邱議瑩"
$legId   = 913

# ---- Header row (row 1): copy the bold/bordered header style from an
# ---- existing header cell (B1) onto the new header cells, then write
# ---- the column names.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item(1, 8).Value  = "date"
$ws.Cells.Item(1, 9).Value  = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# ---- Data rows (rows 2-3): fill the date/name/id for every existing
# ---- record. The date must stay a literal text string (not get
# ---- auto-converted into a date serial number), so it is written via
# ---- a quoted-text formula and then flattened back down to a plain
# ---- value with Paste Values - this keeps the cell as shared-string
# ---- text without requiring any new style/number-format to be
# ---- registered in the workbook.
$lastRow = $ws.Cells.Item(1, 1).CurrentRegion.Rows.Count
if (-not $lastRow -or $lastRow -lt 3) { $lastRow = 3 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Formula = "=""" + $legDate + """"
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues

    $ws.Cells.Item($r, 9).Value  = $legName
    $ws.Cells.Item($r, 10).Value = $legId
}
